# Add San Juan PR to list of photographed cities, and fill in mileage data
# for several other already-photographed cities. Also clears the leftover
# "0" placeholder values left in the miles column for cities that have not
# been photographed yet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Cities")

# --- San Juan PR (row 92): set photo_date and miles -----------------------
$ws.Range("E92").Value = "2025-01"
$ws.Range("F92").Value = 32

# --- Fill in miles for cities that already had a photo_date but were ------
# --- still missing a mileage figure ----------------------------------------
$ws.Range("F54").Value = 2.8
$ws.Range("F55").Value = 3.7
$ws.Range("F64").Value = 4
$ws.Range("F65").Value = 3.1
$ws.Range("F70").Value = 1
$ws.Range("F71").Value = 1.8
$ws.Range("F78").Value = 5.1
$ws.Range("F79").Value = 7.3
$ws.Range("F81").Value = 3.5
$ws.Range("F83").Value = 6.5
$ws.Range("F89").Value = 7.7
$ws.Range("F90").Value = 26.2
$ws.Range("F91").Value = 7.7
$ws.Range("F116").Value = 8.9
$ws.Range("F117").Value = 5.6

# --- Clear the stray "0" placeholders left in the miles column for --------
# --- cities that still have not been photographed --------------------------
$emptyMilesRows = @(2,3,4,7,8,9,10,11,12,15,17,18,19,20,21,34,37,38,43,46,60,61,62,63,76,85,86,87,88,100,102,103,104,105,106,107,113,114,115,120,121)
foreach ($r in $emptyMilesRows) {
    $ws.Range("F$r").ClearContents()
}
